$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First, push the existing placeholder row 7 (style + content) down to rows 8-16 ---
for ($r = 8; $r -le 16; $r++) {
    $ws.Range("H7").Copy($ws.Cells.Item($r, 8))
}

# --- Now fill in the new row 7 data (Best Time to Buy and Sell Stock) ---
$ws.Range("A7").Value = 121
$ws.Range("B7").Value = "Best Time to Buy and Sell Stock"
$ws.Range("C7").Value = "Array"
$ws.Range("D7").Value = "Array, Dynamic Prog"
$ws.Range("F7").Value = "Easy"
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = "✅"
$ws.Range("I7").Value = "Given O(n) sol but not accepted hence sol by chatGPT without understanding"

# A7 gets the blue "done" highlight fill used elsewhere in the sheet (e.g. A5)
$ws.Range("A7").Interior.Color = 12611584

# H7 gets the red checkmark-style font used elsewhere in the sheet (e.g. H2:H6)
$ws.Range("H7").ClearFormats()
$ws.Range("H7").Value = "✅"
$ws.Range("H7").Font.Color = 255

# Row 7 is taller to fit the wrapped text
$ws.Rows.Item(7).RowHeight = 60

# --- Update the active selection to I9 ---
$ws.Range("I9").Select() | Out-Null

Write-Host "done"
